# Project version numbers.xlsx - add a new release column (10.5.0) and a new
# project row (Tardigrade.Framework.ZXingNet), per the commit:
#   "Added a new project that implements a QR Code processor based upon
#   ZXing.Net. Migrated existing xUnit test projects from the Tardigrade
#   Framework Solution to a new Solution. Updated NuGet packages."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column D: release "10.5.0" values for the existing projects, plus
#    the new project's own version (row 11), filled top-to-bottom first.
# ---------------------------------------------------------------------
$ws.Range("D1").Value2 = "10.5.0"
$ws.Range("D2").Value2 = "8.8.0"
$ws.Range("D3").Value2 = "4.1.3"
$ws.Range("D4").Value2 = "3.2.3"
$ws.Range("D5").Value2 = "1.2.1"
$ws.Range("D6").Value2 = "5.3.2"
$ws.Range("D7").Value2 = "9.2.1"
$ws.Range("D8").Value2 = "8.1.2"
$ws.Range("D9").Value2 = "1.2.4"
$ws.Range("D10").Value2 = "3.2.3"
$ws.Range("D11").Value2 = "1.0.0"

# ---------------------------------------------------------------------
# 2. New row 11: a new project, Tardigrade.Framework.ZXingNet, introduced
#    at release 10.5.0 (it has no version for the two earlier releases).
# ---------------------------------------------------------------------
$ws.Range("A11").Value2 = "Tardigrade.Framework.ZXingNet"
$ws.Range("B11").Value2 = "-"
$ws.Range("C11").Value2 = "-"

# ---------------------------------------------------------------------
# 3. Formatting: header row stays bold; the "no version change between
#    releases" cells keep their de-emphasised (non-bold, non-italic)
#    styling, matching the same rows that were already de-emphasised in
#    column C (rows 3,4,6,8,9,10 plus the new row 11).
# ---------------------------------------------------------------------
$ws.Range("A1:D1").Font.Bold = $true

foreach ($r in 3,4,6,8,9,10) {
    $ws.Range("C$r").Font.Italic = $false
    $ws.Range("D$r").Font.Italic = $false
}
$ws.Range("C11").Font.Italic = $false
$ws.Range("D11").Font.Italic = $false

# ---------------------------------------------------------------------
# 4. Selection / view bookkeeping to match the author's last saved state.
# ---------------------------------------------------------------------
$ws.Range("D11").Select() | Out-Null
